$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 4140992
$ws.Range("H2").Value = 676320
$ws.Range("I2").Value = 2.58
$ws.Range("J2").Value = 15.7966
$ws.Range("G3").Value = 3666632
$ws.Range("H3").Value = 744372
$ws.Range("I3").Value = 2.8395999999999999
$ws.Range("J3").Value = 13.9871
$ws.Range("G4").Value = 3525992
$ws.Range("H4").Value = 807801
$ws.Range("I4").Value = 3.0815000000000001
$ws.Range("J4").Value = 13.4506
$ws.Range("G5").Value = 3305336
$ws.Range("H5").Value = 995187
$ws.Range("I5").Value = 3.7963
$ws.Range("J5").Value = 12.6089
$ws.Range("G7").Value = 1268926
$ws.Range("H7").Value = 256089
$ws.Range("I7").Value = 0.97689999999999999
$ws.Range("J7").Value = 4.8406000000000002
$ws.Range("G8").Value = 1230050
$ws.Range("H8").Value = 263685
$ws.Range("I8").Value = 1.0059
$ws.Range("J8").Value = 4.6923000000000004
$ws.Range("G9").Value = 1217088
$ws.Range("H9").Value = 271923
$ws.Range("I9").Value = 1.0373000000000001
$ws.Range("J9").Value = 4.6428000000000003
$ws.Range("G10").Value = 1192690
$ws.Range("H10").Value = 295474
$ws.Range("I10").Value = 1.1271
$ws.Range("J10").Value = 4.5498000000000003
$ws.Range("G12").Value = 19878280
$ws.Range("H12").Value = 1347795
$ws.Range("I12").Value = 5.1414
$ws.Range("J12").Value = 75.829599999999999
$ws.Range("C13").Value = 1551753
$ws.Range("D13").Value = 775877
$ws.Range("G13").Value = 18615496
$ws.Range("H13").Value = 1612530
$ws.Range("I13").Value = 6.1513
$ws.Range("J13").Value = 71.012500000000003
$ws.Range("C14").Value = 1245343
$ws.Range("D14").Value = 622672
$ws.Range("G14").Value = 18113938
$ws.Range("H14").Value = 1908663
$ws.Range("I14").Value = 7.2809999999999997
$ws.Range("J14").Value = 69.099199999999996
$ws.Range("G15").Value = 17417504
$ws.Range("H15").Value = 2544567
$ws.Range("I15").Value = 9.7067999999999994
$ws.Range("J15").Value = 66.442499999999995
$ws.Range("G17").Value = 35167778
$ws.Range("H17").Value = 3180800
$ws.Range("I17").Value = 12.133800000000001
$ws.Range("J17").Value = 134.15440000000001
$ws.Range("G18").Value = 32094392
$ws.Range("H18").Value = 3465943
$ws.Range("I18").Value = 13.221500000000001
$ws.Range("J18").Value = 122.43040000000001
$ws.Range("G19").Value = 31282248
$ws.Range("H19").Value = 3834671
$ws.Range("I19").Value = 14.6281
$ws.Range("J19").Value = 119.3323
$ws.Range("G20").Value = 29753636
$ws.Range("H20").Value = 5020619
$ws.Range("I20").Value = 19.152100000000001
$ws.Range("J20").Value = 113.50109999999999

$ws.Range("J20").Select() | Out-Null
